# Weekly update: insert two new rows of fresh data (488:489) at the top of
# the "Betarraga" price block, pushing the existing 488:529 block down to
# 490:531 (dimension grows from A1:R529 to A1:R531).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right before the current row 488; this shifts all
# rows from 488 downward by 2 (old 488 -> 490, ..., old 529 -> 531) and
# carries the existing row formatting (e.g. the date style on column D)
# into the freshly inserted rows.
$ws.Range("A488:R489").EntireRow.Insert()

# Row 488: new "Primera" quality record
$ws.Range("A488").Value = 3
$ws.Range("B488").Value = "Femacal de La Calera"
$ws.Range("C488").Value = "Coquimbo"
$ws.Range("D488").Value = 44578
$ws.Range("E488").Value = 5
$ws.Range("F488").Value = 100114014
$ws.Range("G488").Value = "Betarraga"
$ws.Range("H488").Value = "Sin especificar"
$ws.Range("I488").Value = "Primera"
$ws.Range("J488").Value = 3400
$ws.Range("K488").Value = 500
$ws.Range("L488").Value = 550
$ws.Range("M488").Value = 526
$ws.Range("N488").Value = "`$/paquete 4 unidades"
$ws.Range("O488").Value = "Provincia de Quillota"
$ws.Range("P488").Value = 132
$ws.Range("Q488").Value = 4
$ws.Range("R488").Value = "Hortaliza"

# Row 489: new "Segunda" quality record
$ws.Range("A489").Value = 3
$ws.Range("B489").Value = "Femacal de La Calera"
$ws.Range("C489").Value = "Coquimbo"
$ws.Range("D489").Value = 44578
$ws.Range("E489").Value = 5
$ws.Range("F489").Value = 100114014
$ws.Range("G489").Value = "Betarraga"
$ws.Range("H489").Value = "Sin especificar"
$ws.Range("I489").Value = "Segunda"
$ws.Range("J489").Value = 1800
$ws.Range("K489").Value = 400
$ws.Range("L489").Value = 400
$ws.Range("M489").Value = 400
$ws.Range("N489").Value = "`$/paquete 4 unidades"
$ws.Range("O489").Value = "Provincia de Quillota"
$ws.Range("P489").Value = 100
$ws.Range("Q489").Value = 4
$ws.Range("R489").Value = "Hortaliza"
